$p = $ppt.ActivePresentation
Write-Output ($p.ColorSchemes | Get-Member | Out-String)
